# Reverse the order of the comma-separated "Recorded By" entries in column G,
# e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System".
#
# The single combination "System, admin@admin.com" is intentionally left
# untouched, matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Text

    if ($null -eq $val -or $val -eq '') {
        continue
    }

    if ($val -eq 'System, admin@admin.com') {
        continue
    }

    $parts = $val -split ', '
    if ($parts.Count -lt 2) {
        continue
    }

    $reversed = @()
    for ($i = $parts.Count - 1; $i -ge 0; $i--) {
        $reversed += $parts[$i]
    }

    $cell.Value = [string]::Join(', ', $reversed)
}
